# Mise à jour de l'application
# Adds 5 new "Entrainement" (training) GPS rows (J+3, 2025-09-02) for
# five players, appended after the existing data (rows 412-416).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (player, position, GPS metrics) -----------------------
# Column layout: A Type | B Date | C Scope | D MatchDay | E Player |
#                F Position | G Temps joué | H..V GPS metrics

$newRows = @(
    @{ Row=412; E="Omar Benyounes";   F="center midfield"; G="01:31:50";
       H=8.2799999999999994; I=0.4;  J=7.87; K=0.41; L=0;    M=0;    N=0;    O=0;
       P=5.32; Q=20.13; R=4.4000000000000004; S=41; T=4;  U=43; V=5 },
    @{ Row=413; E="Karahali Souaré";  F="right forward";    G="01:29:24";
       H=6.73; I=0.32; J=6.4;  K=0.25; L=0.08; M=0.01; N=0;    O=3;
       P=4.21; Q=26.13; R=5.4;                     S=76; T=20; U=50; V=15 },
    @{ Row=414; E="Naim Dhib";        F="center midfield";  G="01:31:51";
       H=7.66; I=0.37; J=7.28; K=0.35; L=0.03; M=0;    N=0;    O=0;
       P=4.93; Q=21.72; R=4.9400000000000004; S=45; T=5;  U=36; V=5 },
    @{ Row=415; E="Jeremie Laurent";  F="left forward";     G="01:30:47";
       H=7.81; I=0.46; J=7.34; K=0.44; L=0.03; M=0;    N=0;    O=0;
       P=5.0999999999999996; Q=21.99; R=4.63; S=41; T=9;  U=40; V=10 },
    @{ Row=416; E="Mattheo Haon";     F="right back";       G="01:33:14";
       H=8.01; I=0.4;  J=7.6;  K=0.37; L=0.04; M=0;    N=0;    O=0;
       P=5.07; Q=22.71; R=4.8499999999999996; S=49; T=6;  U=33; V=6 }
)

# First pass: fill the "Temps joué" column (G) for every new row. This
# mirrors the order new shared strings were introduced when the source
# application last saved the file (time strings before the "J+3" label).
foreach ($r in $newRows) {
    $ws.Range("G" + $r.Row).Value = $r.G
}

# Second pass: the constant columns shared by every new row.
for ($row = 412; $row -le 416; $row++) {
    $ws.Range("A" + $row).Value = "Entrainement"
    $ws.Range("C" + $row).Value = "Global"
    $ws.Range("D" + $row).Value = "J+3"
}

# Third pass: remaining per-row columns (date, player, position, metrics).
foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("B" + $row).Value = 45902
    $ws.Range("E" + $row).Value = $r.E
    $ws.Range("F" + $row).Value = $r.F

    $ws.Range("H" + $row).Value = $r.H
    $ws.Range("I" + $row).Value = $r.I
    $ws.Range("J" + $row).Value = $r.J
    $ws.Range("K" + $row).Value = $r.K
    $ws.Range("L" + $row).Value = $r.L
    $ws.Range("M" + $row).Value = $r.M
    $ws.Range("N" + $row).Value = $r.N
    $ws.Range("O" + $row).Value = $r.O
    $ws.Range("P" + $row).Value = $r.P
    $ws.Range("Q" + $row).Value = $r.Q
    $ws.Range("R" + $row).Value = $r.R
    $ws.Range("S" + $row).Value = $r.S
    $ws.Range("T" + $row).Value = $r.T
    $ws.Range("U" + $row).Value = $r.U
    $ws.Range("V" + $row).Value = $r.V
}

# --- Formatting: reuse the existing number-format styles instead of ------
# letting Excel mint brand-new ones (date style on column B, centred style
# on column D), exactly like the rows immediately above them.
$ws.Range("B411").Copy() | Out-Null
$ws.Range("B412:B416").PasteSpecial(-4122) | Out-Null

$ws.Range("D397").Copy() | Out-Null
$ws.Range("D412:D416").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- View state: scrolled down & selection on the new last data row ------
$ws.Range("F419").Select() | Out-Null
